# Trade #74 closed at 2026-02-18 00:27:44 - unknown UNKNOWN +0.000%
#
# This script:
#  1. Updates the Summary sheet roll-up metrics.
#  2. Updates the Strategy Status sheet's MarketMaking row.
#  3. Closes trade #102 (row 103 on "All Trades", row 35 on "MarketMaking")
#     which now shows an early exit with a small profit.
#  4. Appends a new open trade #131 (row 132 on "All Trades", row 52 on
#     "MarketMaking").

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, [string]$text) {
    # Force a literal string even when the text looks like a date/time so
    # Excel's COM layer doesn't silently convert it to a date serial.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.33   # Current Capital
$summary.Range("B4").Value = 0.44      # Total P&L $
$summary.Range("B5").Value = 0.09      # Total P&L %
$summary.Range("B6").Value = 102       # Total Trades
$summary.Range("B7").Value = 48        # Winning Trades
$summary.Range("B9").Value = 47.06     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.48
$status.Range("D6").Value = 34
$status.Range("E6").Value = -0.33
$status.Range("F6").Value = -0.52
$status.Range("G6").Value = 47.06

# ---------------------------------------------------------------------
# 3) All Trades sheet - close trade #102 (row 103) + append trade #131
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G103").Value = 0.77
$allTrades.Range("H103").Value = "CLOSED"
$allTrades.Range("I103").Value = 16.6667
$allTrades.Range("J103").Value = 0.11
$allTrades.Range("K103").Value = 99.48
Set-TextValue $allTrades.Range("L103") "early_exit"
$allTrades.Range("M103").Value = 0.13

$allTrades.Cells.Item(132, 1).Value = 131
Set-TextValue $allTrades.Cells.Item(132, 2) "2026-02-18"
Set-TextValue $allTrades.Cells.Item(132, 3) "00:27:39"
Set-TextValue $allTrades.Cells.Item(132, 4) "MarketMaking"
Set-TextValue $allTrades.Cells.Item(132, 5) "DOWN"
$allTrades.Cells.Item(132, 6).Value = 0.66
$allTrades.Cells.Item(132, 8).Value = "OPEN"
$allTrades.Cells.Item(132, 9).Value = 0
$allTrades.Cells.Item(132, 10).Value = 0
$allTrades.Cells.Item(132, 11).Value = 99.36967800952272
$allTrades.Cells.Item(132, 13).Value = 0
$allTrades.Cells.Item(132, 14).Value = 0
$allTrades.Cells.Item(132, 15).Value = 0
$allTrades.Cells.Item(132, 16).Value = 0.65
Set-TextValue $allTrades.Cells.Item(132, 17) "Wide spread capture: 392 bps vs avg 295 bps"

# ---------------------------------------------------------------------
# 4) MarketMaking sheet - close trade #102 (row 35) + append trade #131
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Range("G35").Value = 0.77
$mm.Range("H35").Value = "CLOSED"
$mm.Range("I35").Value = 16.6667
$mm.Range("J35").Value = 0.11
$mm.Range("K35").Value = 99.48
Set-TextValue $mm.Range("P35") "early_exit"
$mm.Range("Q35").Value = 0.13

$mm.Cells.Item(52, 1).Value = 131
Set-TextValue $mm.Cells.Item(52, 2) "2026-02-18"
Set-TextValue $mm.Cells.Item(52, 3) "00:27:39"
Set-TextValue $mm.Cells.Item(52, 4) "MarketMaking"
Set-TextValue $mm.Cells.Item(52, 5) "DOWN"
$mm.Cells.Item(52, 6).Value = 0.66
$mm.Cells.Item(52, 8).Value = "OPEN"
$mm.Cells.Item(52, 9).Value = 0
$mm.Cells.Item(52, 10).Value = 0
$mm.Cells.Item(52, 11).Value = 99.36967800952272
$mm.Cells.Item(52, 12).Value = 0
$mm.Cells.Item(52, 13).Value = 0
$mm.Cells.Item(52, 14).Value = 0.65
Set-TextValue $mm.Cells.Item(52, 15) "Wide spread capture: 392 bps vs avg 295 bps"
$mm.Cells.Item(52, 17).Value = 0
